$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "ChEBI_DEF" header, matching the header style used by B1:E1
$ws.Range("F1").Value = "ChEBI_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("F1").VerticalAlignment = -4160    # xlTop
$ws.Range("F1").Borders.LineStyle = 1

# New column F value for row 2
$ws.Range("F2").Value = "[]"
